# "Update countries & provincias Spain" -- refresh COVID country/region stats
# (new timestamp, a handful of relabeled countries caused by a shared-string
# resort, and updated case counts for many countries).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A text relabels (caused by shared-string reordering) ---
$ws.Range("A1").Value = 'Datos actualizados a 5 de Mayo de 2020 a las 21:03'
$ws.Range("A16").Value = 'Peru'
$ws.Range("A17").Value = 'Belgica'
$ws.Range("A81").Value = 'Cuba'
$ws.Range("A82").Value = 'Bolivia'
$ws.Range("A197").Value = 'Curazao'
$ws.Range("A198").Value = 'Dominica'
$ws.Range("A199").Value = 'San Cristobal y Nieves'
$ws.Range("A200").Value = 'Burundi'

# --- Updated numeric data (columns B-H) ---
$ws.Range("B4").Value = 1226123
$ws.Range("C4").Value = 13288
$ws.Range("D4").Value = 192256
$ws.Range("E4").Value = 962563
$ws.Range("G4").Value = 1383
$ws.Range("H4").Value = 71304
$ws.Range("B8").Value = 170551
$ws.Range("C8").Value = 1089
$ws.Range("D8").Value = 52736
$ws.Range("E8").Value = 92284
$ws.Range("F8").Value = 3430
$ws.Range("G8").Value = 330
$ws.Range("H8").Value = 25531
$ws.Range("B9").Value = 166521
$ws.Range("C9").Value = 369
$ws.Range("E9").Value = 24428
$ws.Range("B12").Value = 110156
$ws.Range("C12").Value = 1890
$ws.Range("E12").Value = 56856
$ws.Range("G12").Value = 142
$ws.Range("H12").Value = 7485
$ws.Range("B15").Value = 61959
$ws.Range("C15").Value = 1187
$ws.Range("D15").Value = 26650
$ws.Range("E15").Value = 31273
$ws.Range("G15").Value = 182
$ws.Range("H15").Value = 4036
$ws.Range("B16").Value = 51189
$ws.Range("C16").Value = 3817
$ws.Range("D16").Value = 14427
$ws.Range("E16").Value = 35318
$ws.Range("F16").Value = 709
$ws.Range("G16").Value = 100
$ws.Range("H16").Value = 1444
$ws.Range("B17").Value = 50509
$ws.Range("C17").Value = 242
$ws.Range("D17").Value = 12441
$ws.Range("E17").Value = 30052
$ws.Range("F17").Value = 646
$ws.Range("G17").Value = 92
$ws.Range("H17").Value = 8016
$ws.Range("B18").Value = 49400
$ws.Range("C18").Value = 2963
$ws.Range("D18").Value = 14142
$ws.Range("E18").Value = 33565
$ws.Range("G18").Value = 127
$ws.Range("H18").Value = 1693
$ws.Range("D22").Value = 25400
$ws.Range("E22").Value = 2814
$ws.Range("B32").Value = 16289
$ws.Range("C32").Value = 43
$ws.Range("D32").Value = 10465
$ws.Range("E32").Value = 5586
$ws.Range("G32").Value = 3
$ws.Range("H32").Value = 238
$ws.Range("B72").Value = 2207
$ws.Range("C72").Value = 18
$ws.Range("E72").Value = 696
$ws.Range("B81").Value = 1685
$ws.Range("C81").Value = 17
$ws.Range("D81").Value = 954
$ws.Range("E81").Value = 662
$ws.Range("F81").Value = 8
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 69
$ws.Range("B82").Value = 1681
$ws.Range("C82").Value = 87
$ws.Range("D82").Value = 174
$ws.Range("E82").Value = 1425
$ws.Range("F82").Value = 3
$ws.Range("G82").Value = 6
$ws.Range("H82").Value = 82
$ws.Range("B100").Value = 771
$ws.Range("C100").Value = 20
$ws.Range("E100").Value = 549
$ws.Range("B115").Value = 573
$ws.Range("C115").Value = 32
$ws.Range("E115").Value = 555
$ws.Range("B129").Value = 326
$ws.Range("C129").Value = 1
$ws.Range("E129").Value = 32
$ws.Range("F129").Value = 19
$ws.Range("D159").Value = 81
$ws.Range("E159").Value = 10
$ws.Range("D177").Value = 30
$ws.Range("E177").Value = 9
$ws.Range("B179").Value = 36
$ws.Range("C179").Value = 1
$ws.Range("E179").Value = 23
$ws.Range("E197").Value = 2
$ws.Range("H197").Value = 1
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 0
$ws.Range("D199").Value = 8
$ws.Range("H199").Value = 0
$ws.Range("D200").Value = 7
$ws.Range("H200").Value = 1
